# Add two new asset-class columns ("Chinese Domestic Equity" and "Hong Kong
# Equity") to the scenario table on the "Value" sheet. Both new columns
# mirror the existing "AC World Equity" (column D) scenario shocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "Chinese Domestic Equity"
$ws.Range("J1").Value = "Hong Kong Equity"

# Data rows 2-8: same shock values as column D (AC World Equity)
$values = @(
    -0.20532831707423782,
    -0.48970981536305258,
    -0.59066781585282491,
    -0.34102902645985778,
    -0.265932178785099,
    -0.16786465884711477,
    -0.1
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $values[$i]
    $ws.Cells.Item($row, 10).Value = $values[$i]
}

# Restore the active selection to A2 (matches saved view state in the file)
$ws.Range("A2").Select()
